$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.558564
$ws.Range("H2").Value = 4.675692
$ws.Range("I2").Value = 0.005692101168584756
$ws.Range("J2").Value = 0.005692101168584756
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.131233999999999
$ws.Range("N2").Value = 24.393702
$ws.Range("O2").Value = 0.02090995573015822
$ws.Range("P2").Value = 0.02090995573015823
$ws.Range("Q2").Value = 12.673048587976
$ws.Range("R2").Value = 114.057437291784
$ws.Range("S2").Value = 0.0001190215834466891
$ws.Range("T2").Value = 0.0001190215834466892

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.558564
$ws.Range("H3").Value = 4.675692
$ws.Range("I3").Value = 0.005692101168584756
$ws.Range("J3").Value = 0.005692101168584756
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 243.3763986666667
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.625857000534647
$ws.Range("P3").Value = 0.6258570005346471
$ws.Range("Q3").Value = 379.3176934115147
$ws.Range("R3").Value = 3413.859240703632
$ws.Range("S3").Value = 0.003562441364110215
$ws.Range("T3").Value = 0.003562441364110215

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.558564
$ws.Range("H4").Value = 4.675692
$ws.Range("I4").Value = 0.005692101168584756
$ws.Range("J4").Value = 0.005692101168584756
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 103.9426383333333
$ws.Range("N4").Value = 311.827915
$ws.Range("O4").Value = 0.2672947262403034
$ws.Range("P4").Value = 0.2672947262403035
$ws.Range("Q4").Value = 162.0012541713533
$ws.Range("R4").Value = 1458.01128754218
$ws.Range("S4").Value = 0.001521468623588974
$ws.Range("T4").Value = 0.001521468623588974

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.558564
$ws.Range("H5").Value = 4.675692
$ws.Range("I5").Value = 0.005692101168584756
$ws.Range("J5").Value = 0.005692101168584756
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.41874933333333
$ws.Range("N5").Value = 100.256248
$ws.Range("O5").Value = 0.08593831749489127
$ws.Range("P5").Value = 0.08593831749489128
$ws.Range("Q5").Value = 52.08525963595732
$ws.Range("R5").Value = 468.767336723616
$ws.Range("S5").Value = 0.0004891695974388783
$ws.Range("T5").Value = 0.0004891695974388784

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 259.5505726666667
$ws.Range("H6").Value = 778.6517180000001
$ws.Range("I6").Value = 0.9479162344201305
$ws.Range("J6").Value = 0.9479162344201304
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.131233999999999
$ws.Range("N6").Value = 24.393702
$ws.Range("O6").Value = 0.02090995573015822
$ws.Range("P6").Value = 0.02090995573015823
$ws.Range("Q6").Value = 2110.466441186671
$ws.Range("R6").Value = 18994.19797068004
$ws.Range("S6").Value = 0.01982088649762321
$ws.Range("T6").Value = 0.01982088649762321

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 259.5505726666667
$ws.Range("H7").Value = 778.6517180000001
$ws.Range("I7").Value = 0.9479162344201305
$ws.Range("J7").Value = 0.9479162344201304
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 243.3763986666667
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.625857000534647
$ws.Range("P7").Value = 0.6258570005346471
$ws.Range("Q7").Value = 63168.48364748432
$ws.Range("R7").Value = 568516.3528273589
$ws.Range("S7").Value = 0.5932600112322802
$ws.Range("T7").Value = 0.5932600112322802

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 259.5505726666667
$ws.Range("H8").Value = 778.6517180000001
$ws.Range("I8").Value = 0.9479162344201305
$ws.Range("J8").Value = 0.9479162344201304
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 103.9426383333333
$ws.Range("N8").Value = 311.827915
$ws.Range("O8").Value = 0.2672947262403034
$ws.Range("P8").Value = 0.2672947262403035
$ws.Range("Q8").Value = 26978.37130390089
$ws.Range("R8").Value = 242805.341735108
$ws.Range("S8").Value = 0.2533730103780681
$ws.Range("T8").Value = 0.2533730103780681

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 259.5505726666667
$ws.Range("H9").Value = 778.6517180000001
$ws.Range("I9").Value = 0.9479162344201305
$ws.Range("J9").Value = 0.9479162344201304
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 33.41874933333333
$ws.Range("N9").Value = 100.256248
$ws.Range("O9").Value = 0.08593831749489127
$ws.Range("P9").Value = 0.08593831749489128
$ws.Range("Q9").Value = 8673.855527270453
$ws.Range("R9").Value = 78064.69974543407
$ws.Range("S9").Value = 0.08146232631215895
$ws.Range("T9").Value = 0.08146232631215895

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 12.18925266666667
$ws.Range("H10").Value = 36.567758
$ws.Range("I10").Value = 0.04451691386950307
$ws.Range("J10").Value = 0.04451691386950307
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.131233999999999
$ws.Range("N10").Value = 24.393702
$ws.Range("O10").Value = 0.02090995573015822
$ws.Range("P10").Value = 0.02090995573015823
$ws.Range("Q10").Value = 99.11366571779065
$ws.Range("R10").Value = 892.0229914601158
$ws.Range("S10").Value = 0.0009308466982545758
$ws.Range("T10").Value = 0.0009308466982545759

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 12.18925266666667
$ws.Range("H11").Value = 36.567758
$ws.Range("I11").Value = 0.04451691386950307
$ws.Range("J11").Value = 0.04451691386950307
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 243.3763986666667
$ws.Range("N11").Value = 730.1291960000001
$ws.Range("O11").Value = 0.625857000534647
$ws.Range("P11").Value = 0.6258570005346471
$ws.Range("Q11").Value = 2966.576416451397
$ws.Range("R11").Value = 26699.18774806257
$ws.Range("S11").Value = 0.02786122218742642
$ws.Range("T11").Value = 0.02786122218742642

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 12.18925266666667
$ws.Range("H12").Value = 36.567758
$ws.Range("I12").Value = 0.04451691386950307
$ws.Range("J12").Value = 0.04451691386950307
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 103.9426383333333
$ws.Range("N12").Value = 311.827915
$ws.Range("O12").Value = 0.2672947262403034
$ws.Range("P12").Value = 0.2672947262403035
$ws.Range("Q12").Value = 1266.983081484952
$ws.Range("R12").Value = 11402.84773336457
$ws.Range("S12").Value = 0.01189913630581199
$ws.Range("T12").Value = 0.01189913630581199

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 12.18925266666667
$ws.Range("H13").Value = 36.567758
$ws.Range("I13").Value = 0.04451691386950307
$ws.Range("J13").Value = 0.04451691386950307
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 33.41874933333333
$ws.Range("N13").Value = 100.256248
$ws.Range("O13").Value = 0.08593831749489127
$ws.Range("P13").Value = 0.08593831749489128
$ws.Range("Q13").Value = 407.3495794279982
$ws.Range("R13").Value = 3666.146214851984
$ws.Range("S13").Value = 0.003825708678010084
$ws.Range("T13").Value = 0.003825708678010084

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.5133286666666667
$ws.Range("H14").Value = 1.539986
$ws.Range("I14").Value = 0.001874750541781658
$ws.Range("J14").Value = 0.001874750541781658
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 8.131233999999999
$ws.Range("N14").Value = 24.393702
$ws.Range("O14").Value = 0.02090995573015822
$ws.Range("P14").Value = 0.02090995573015823
$ws.Range("Q14").Value = 4.173995507574666
$ws.Range("R14").Value = 37.56595956817199
$ws.Range("S14").Value = 0.00003920095083374462
$ws.Range("T14").Value = 0.00003920095083374462

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.5133286666666667
$ws.Range("H15").Value = 1.539986
$ws.Range("I15").Value = 0.001874750541781658
$ws.Range("J15").Value = 0.001874750541781658
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 243.3763986666667
$ws.Range("N15").Value = 730.1291960000001
$ws.Range("O15").Value = 0.625857000534647
$ws.Range("P15").Value = 0.6258570005346471
$ws.Range("Q15").Value = 124.9320822256951
$ws.Range("R15").Value = 1124.388740031256
$ws.Range("S15").Value = 0.001173325750830173
$ws.Range("T15").Value = 0.001173325750830173

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.5133286666666667
$ws.Range("H16").Value = 1.539986
$ws.Range("I16").Value = 0.001874750541781658
$ws.Range("J16").Value = 0.001874750541781658
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 103.9426383333333
$ws.Range("N16").Value = 311.827915
$ws.Range("O16").Value = 0.2672947262403034
$ws.Range("P16").Value = 0.2672947262403035
$ws.Range("Q16").Value = 53.35673594546555
$ws.Range("R16").Value = 480.21062350919
$ws.Range("S16").Value = 0.0005011109328343889
$ws.Range("T16").Value = 0.0005011109328343889

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.5133286666666667
$ws.Range("H17").Value = 1.539986
$ws.Range("I17").Value = 0.001874750541781658
$ws.Range("J17").Value = 0.001874750541781658
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 33.41874933333333
$ws.Range("N17").Value = 100.256248
$ws.Range("O17").Value = 0.08593831749489127
$ws.Range("P17").Value = 0.08593831749489128
$ws.Range("Q17").Value = 17.15480203694755
$ws.Range("R17").Value = 154.393218332528
$ws.Range("S17").Value = 0.0001611129072833515
$ws.Range("T17").Value = 0.0001611129072833515
